$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append row 5 with the results of the 2025-11-29 run
$ws.Range("A5").Value = "'11/29/2025"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = 14475.31
$ws.Range("C5").Value = 0.1631502794106506
$ws.Range("D5").Value = 0.8368497205893494
$ws.Range("E5").Value = -56.34
$ws.Range("F5").Value = -13.74
$ws.Range("G5").Value = -18353.1
$ws.Range("H5").Value = -60.24
$ws.Range("I5").Value = -391.2
$ws.Range("J5").Value = -14.21
